$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (NIFTY ...18050 CE -> ...17750 CE) plus related numeric-as-text fields
$ws.Range("J2").Value = "NIFTY2240717750CE"
$ws.Range("P2").Value = "'3"
$ws.Range("R2").Value = "'1"
$ws.Range("T2").Value = "'5"
$ws.Range("Z2").Value = "'0.6"

# Row 3 (NIFTY ...18050 PE -> ...17750 PE) plus related numeric-as-text fields
$ws.Range("J3").Value = "NIFTY2240717750PE"
$ws.Range("P3").Value = "'3"
$ws.Range("R3").Value = "'1"
$ws.Range("T3").Value = "'5"
$ws.Range("Z3").Value = "'0.6"
